# Regenerate the "K" (strikeouts) column (G) values in save_data for edwards_carl.
# The data source for this column was regenerated (switching from a prior
# "Strike#" derived value to the actual K count), so the new values below are
# written directly as the refreshed source data, matching the regen output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 0
    36 = 3
    37 = 1
    38 = 1
    39 = 2
    40 = 1
    41 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
